# SSDM-55: fixed xls export types data.
# Adds a new "Multivalued" column (K) to the data-set-type export sample sheet,
# with a bold header and "TRUE"/"FALSE" formatted boolean-looking text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New boolean-style display format used for the Multivalued column values
# (shows the word TRUE/FALSE while the underlying cell keeps the literal
# text "FALSE", matching the rest of the boolean-like columns on the sheet).
$boolFormat = """TRUE"";""TRUE"";""FALSE"""

# Header cell K4: "Multivalued", bold black Calibri 11 (same look as the
# other header cells in row 4, just bold).
$ws.Range("K4").Value = "Multivalued"
$ws.Range("K4").Font.Bold = $true
$ws.Range("K4").Font.Size = 11

# Data cells K5:K8: text "FALSE" (the leading apostrophe forces the value to
# stay text instead of being auto-converted to a boolean), left aligned, with
# the custom TRUE/FALSE display format. K8 is an extra formatted row with no
# content, matching the other Multivalued cells' formatting. The format and
# alignment are applied before the value so all four cells end up sharing a
# single, identical cell style.
$ws.Range("K5:K8").NumberFormat = $boolFormat
$ws.Range("K5:K8").HorizontalAlignment = -4131
$ws.Range("K5:K8").Value = "'FALSE"
$ws.Range("K8").ClearContents()

# Restyle the row-3 cells so they explicitly carry the (same) default font,
# matching the refreshed export template.
$ws.Range("D3:F3").Font.Name = "Calibri"

# Update the selection to match the newly added column, like the template
# regeneration that produced this sheet.
$ws.Range("K4:K7").Select() | Out-Null
